$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.561.84"
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("D3").Value = "2.901.10"
$ws.Range("E3").Value = "  -2.52%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'525.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").Value = "'142.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.09%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.547"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.81%  "
$ws.Range("D9").Value = "2.908.17"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("E10").Value = "  -5.27%  "
$ws.Range("D11").Value = "'5.97"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.62%  "
$ws.Range("D12").Value = "'0.358"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("D13").Value = "3.409.42"
$ws.Range("E13").Value = "  -2.49%  "
$ws.Range("D14").Value = "'0.129"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.62%  "
$ws.Range("D15").Value = "60.562.75"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "'22.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.97%  "
$ws.Range("D17").Value = "2.914.74"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").Value = "'0.0000140"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.11%  "
$ws.Range("D19").Value = "'4.96"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.83%  "
$ws.Range("D20").Value = "'11.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.06%  "
$ws.Range("D21").Value = "'350.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.17%  "
$ws.Range("D22").Value = "'6.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.78%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'5.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("E25").Value = "  -1.86%  "
$ws.Range("D26").Value = "'0.451"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.28%  "
$ws.Range("D27").Value = "'0.178"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.32%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "'7.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.98%  "
$ws.Range("D30").Value = "0.0₃0848"
$ws.Range("E30").Value = "  -9.64%  "
$ws.Range("E32").Value = "  -2.99%  "
$ws.Range("D33").Value = "'19.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.69%  "
$ws.Range("D34").Value = "'151.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.73%  "
$ws.Range("D35").Value = "'4.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.32%  "
$ws.Range("D36").Value = "'5.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.91%  "
$ws.Range("D37").Value = "'0.994"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.21%  "
$ws.Range("D38").Value = "'1.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.12%  "
$ws.Range("D39").Value = "'37.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "'1.46"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.36%  "
$ws.Range("D41").Value = "'3.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.22%  "
$ws.Range("D42").Value = "2.286.83"
$ws.Range("E42").Value = "  -5.13%  "
$ws.Range("E43").Value = "  -3.42%  "
$ws.Range("D44").Value = "'0.0578"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("D45").Value = "'20.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.01%  "
$ws.Range("D46").Value = "'0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "'4.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.33%  "
$ws.Range("E48").Value = "  -3.34%  "
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("D50").Value = "'0.0915"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.99%  "
$ws.Range("D51").Value = "'18.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.90%  "
